$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target value is written with a leading apostrophe to force
# text interpretation (prevents Excel from coercing numeric-looking
# strings like "583.20" or "0.0440" into floats and dropping the
# trailing zeros / thousands-style dots), then the style is reset
# to Normal so no stray quote-prefix formatting/style index lingers.

$ws.Range("D2").Value = "'68.834.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.00%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.489.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.42%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'583.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.08%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'190.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.84%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.470.07"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -3.61%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.603"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -3.81%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.04%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.202"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -3.82%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.614"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -5.20%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'51.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.72%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000285"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -6.06%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'9.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -5.19%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.015.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.13%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'635.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +5.48%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'68.879.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.13%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'Uniswap"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'12.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.40%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'WrappedEther"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'3.476.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.90%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -2.29%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'18.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.88%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.946"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -5.54%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'17.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -4.45%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'5.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +3.38%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'99.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -3.48%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'4.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -6.60%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -4.41%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'6.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.19%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'10.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -5.26%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'9.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -5.34%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'32.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.75%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'6.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -7.77%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -12.57%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'11.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -5.81%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.109"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -7.17%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'60.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.84%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.713.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -5.62%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.997"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.35%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.0₃0794"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -9.98%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'3.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.37%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'501.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -5.77%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.76%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.367"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -5.66%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.133"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.90%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'34.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -7.07%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0440"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.60%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -7.07%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.05%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.134"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -4.65%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.997"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.60%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'8.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.57%  "
$ws.Range("E51").Style = "Normal"
